# Weekly data refresh: a new price record (fecha 2021-09-16) is published for
# this market/variety, so a new row is inserted at the top of the data block
# (row 84) and every existing record below it shifts down by one row -
# growing the used range from A1:R121 to A1:R122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 84; Excel shifts rows 84:121 down to 85:122 and the
# sheet's dimension grows to R122 automatically.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A84").Value = 7
$ws.Range("B84").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C84").Value = 'Ñuble'
$ws.Range("D84").Value = 44455
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = 100112006
$ws.Range("G84").Value = 'Repollo'
$ws.Range("H84").Value = 'Crespo record'
$ws.Range("I84").Value = 'Primera'
$ws.Range("J84").Value = 600
$ws.Range("K84").Value = 700
$ws.Range("L84").Value = 750
$ws.Range("M84").Value = 725
$ws.Range("N84").Value = '$/unidad'
$ws.Range("O84").Value = 'Provincia de Diguillín'
$ws.Range("P84").Value = 725
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = 'Hortaliza'
